# Sync attendance_reports: swap the order of names in the "Recorded By"
# column (G) wherever it currently reads "dnasr281@gmail.com, System" so
# that it reads "System, dnasr281@gmail.com" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
